$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("attributes")

# Replace the old magma-style "visible" expression with the new molgenis
# expression syntax for the "age" attribute row.
$ws.Range("G4").Value = "{age} >= 18"

# Reflect the last-active cell selection on the attributes sheet.
$ws.Range("G4").Select()
